$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 44, shifting the existing rows 44-88 down to 47-91.
$ws.Rows("44:46").Insert()

# Populate the 3 newly inserted rows with their new data.
$rows = @(
    @{ Row = 44; D = 44460; K = "Cultivar IV Región"; L = "Primera"; M = 175;  N = 20000; O = 20000; P = 20000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 2000; T = 10 },
    @{ Row = 45; D = 44460; K = "Cultivar IV Región"; L = "Segunda"; M = 70;   N = 18000; O = 18000; P = 18000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 1800; T = 10 },
    @{ Row = 46; D = 44460; K = "Cultivar IV Región"; L = "Tercera"; M = 50;   N = 15000; O = 15000; P = 15000; Q = "`$/bandeja 10 kilos"; R = "Provincia de Limarí"; S = 1500; T = 10 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Cells.Item($i, 1).Value = 6
    $ws.Cells.Item($i, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
    $ws.Cells.Item($i, 3).Value = "Metropolitana"
    $ws.Cells.Item($i, 4).Value = $r.D
    $ws.Cells.Item($i, 5).Value = 13
    $ws.Cells.Item($i, 6).Value = "Fruta"
    $ws.Cells.Item($i, 7).Value = 100107
    $ws.Cells.Item($i, 8).Value = "Otros"
    $ws.Cells.Item($i, 9).Value = 100107002
    $ws.Cells.Item($i, 10).Value = "Chirimoya"
    $ws.Cells.Item($i, 11).Value = $r.K
    $ws.Cells.Item($i, 12).Value = $r.L
    $ws.Cells.Item($i, 13).Value = $r.M
    $ws.Cells.Item($i, 14).Value = $r.N
    $ws.Cells.Item($i, 15).Value = $r.O
    $ws.Cells.Item($i, 16).Value = $r.P
    $ws.Cells.Item($i, 17).Value = $r.Q
    $ws.Cells.Item($i, 18).Value = $r.R
    $ws.Cells.Item($i, 19).Value = $r.S
    $ws.Cells.Item($i, 20).Value = $r.T
}
